# "Generate Report for Handback" - update the localization-status report to
# reflect that ed5744f3-2de6-4a75-ad91-39aeb19c2703.md has now been handed
# back (for both the zh-cn and de-de locales).

$wb = $excel.ActiveWorkbook

$statusHandedBack = "Handed back: in sync with en-US"

# --- Overview sheet: update the status cells for the handed-back file ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = $statusHandedBack
$overview.Range("C3").Value = $statusHandedBack

# --- zh-cn sheet: update status + record the handback datetime ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("B3").Value = $statusHandedBack
$zhcn.Range("G3").Value = "2016-03-08 16:47:46"

# --- de-de sheet: update status + record the handback datetime ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("B3").Value = $statusHandedBack
$dede.Range("G3").Value = "2016-03-08 16:48:01"
